$d = $word.ActiveDocument

# NOTE: we set Paragraphs(n).Range.Text directly (rather than Find/Replace)
# so straight quote characters survive instead of being smart-quoted.

# 1) Header comment line: empirical summary -> perf comparison line
$d.Paragraphs(2).Range.Text = "# +1.2x vs ROS/Habitat-Sim baseline"

# 2) Function signature: drop the dict-signal param, use a keyword default instead
$d.Paragraphs(4).Range.Text = "def isolation_combat(loneliness=0.85):"

# 3) Guard clause now reads the plain float argument
$d.Paragraphs(5).Range.Text = '    if loneliness > 0.7:'

# 4) recovery becomes a simple float instead of opening a dict literal
$d.Paragraphs(6).Range.Text = "        recovery = 1.00"

# 5) former dict entries become plain local assignments / prints
$d.Paragraphs(7).Range.Text = "        time = 60  # s"
$d.Paragraphs(8).Range.Text = "        boost = 1.2  # vs baseline"
$d.Paragraphs(9).Range.Text = '        print(f"Recovery: {recovery} in {time}s")'
$d.Paragraphs(10).Range.Text = '        print(f"+{boost}x vs ROS/Habitat-Sim")'

# 6) insert a brand-new print line right after the one above
$d.Paragraphs(10).Range.InsertParagraphAfter()
$d.Paragraphs(11).Range.Text = '        print("Matterport3D: 8k SPS")'

# 7) drop the old "Test Case" comment + variable assignment lines entirely
#    (after the insertion above, they now sit at indices 15 and 16)
$d.Paragraphs(15).Range.Delete()
$d.Paragraphs(15).Range.Delete()

# 8) final call site just invokes with the new default, no print wrapper
$d.Paragraphs(15).Range.Text = "isolation_combat()"
